# Regional steel production data - fuels sheet update
# Fills in missing fuel data rows (JP coking coal, JP steam coal, JP coke,
# RU natural gas) and adds new fuel entries (US coking coal/coke/natural
# gas, IPCC JP coking coal/PCI coal/steam coal/coke, IPCC coking coal)
# along with source-note comments on the emission-factor formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sciencedirectUrl = "https://pdf.sciencedirectassets.com/271097/1-s2.0-S0301421500X01073/1-s2.0-S0301421501001434/main.pdf?x-amz-security-token=AgoJb3JpZ2luX2VjEC8aCXVzLWVhc3QtMSJHMEUCIDycNgcmzDGaem%2B8hqymQUr9KTZdItafx2%2BhltFWNPC6AiEA4sAQrcZ%2F%2BuQ1csgrpfNBYG%2Bv7vHS4HdcfivGxvIcRtoq2gMIGBACGgwwNTkwMDM1NDY4NjUiDOfIyLb7NBWK%2BYzk%2FCq3AxVVI7%2BCKCw9sa5ecoz%2BySXr9i98aqWRO1jaGXHKMhua1AReJbnM5QHrl3x0rsxfvJ5wE6r%2BRUDYrIxtHemLlaEXLTcedV0N1KWL4A2AXt8mdh0ye4n3uryLo5sTaC9Ppzeq%2Bt3nUGeE%2BE4r38w1vVU00uUJJpQYzfi5rv6RItic7nJoHOjtoUW6AB%2FYWqosmp0FqiWALGAu0NK10ThnZpVHTjbQBuorP9cLkYbiZ2PSTqRDi6kCaa0s0yLp8i1N2U5%2BTysZWRJhzUiHyVr4x2p9R7fgBIYkPId8kxlHksd%2BSz2NfVJLgu8Trmm3ctiwC9dLUinhc%2BuLP8BR7EmqTbcFKM6Lml12FSJf6599qox7%2Bxwi7bR5aFyRdcNblArY60oACk9Jc1XXgc105jEj20ynw6AqLrnW73zYyvOtIgkEA5cqJtiElididWTdV4ssEWi8olNymaJ%2BiL4Ei5v03F5rZX9FejOpu6zldO6WQDmtn7EgIjDBSg76i26MDl2UeGqqHocTNIiPupov%2F19v7F7SlT0L2gxKr5heV1rfN3kq1IeIzZVydFjQ4gcXMLEjhQwOkLPWhAswm6qM5gU6tAEwRu29%2F4FA70UDx9ZjIJaP8uSZIjem2QFJ3b%2BZeytYwTYtofqvJIcjWr8x8Qv85Hta66rX%2F9gwkekMoIu7NzWsqsHETFhcZM%2Fp0EFphKn0T5f2NYaNlKQZzXoMEUuyd6AO%2Fkcw3TaOY0WlmvvqDckwz7R2EO02dS0BcxKJkeZ7R%2FN3ykb1i16dUtgGFczv3WmC2BakbTeAK5pw2qao8KGbeQJLiyZ4Z7E8mOc2JNYHNcxy3W0%3D&AWSAccessKeyId=ASIAQ3PHCVTYZDSYP5FM&Expires=1556291036&Signature=B3LotBMPDQQ5Z15FRvqdGSeqgsk%3D&hash=66788d42786115b4239d0f98734c66c7709e09b612f6b0c0a8211ff0d1ae8f55&host=68042c943591013ac2b2430a89b270f6af2c76d8dfd086a07176afe7c76c2c61&pii=S0301421501001434&tid=spdf-9a5bd69d-60c0-4a8b-933f-025252cf964f&sid=a42536342417a14d688b7f856217c6c10482gxrqb&type=client"
$ipccStationarySource = "IPCC EFDB for CO2/TJ: standard for stationary combution in manufacturing indusries and construction"
$ipccRietiSource = "IPCC EFDB source: https://www.rieti.go.jp/users/kainou-kazunari/14j047_e.pdf"
$ipccDensitySource = "IPCC EFDB, density calculated using https://www.unitrove.com/engineering/tools/gas/natural-gas-density"
$noteText = "S.E. Tanzer:`ngCO2/gC * gC/MJ * MJ/kg * g/kg"

# --- Fuel-name labels, entered in sheet-building order so new shared --------
# --- strings land at the same indices the source workbook used. ------------
$ws.Range("A27").Value = "US coking coal"
$ws.Range("A28").Value = "US coke"
$ws.Range("A29").Value = "US natural gas"

$ws.Range("Q20").Value = $sciencedirectUrl
$ws.Range("Q21").Value = $sciencedirectUrl
$ws.Range("Q23").Value = $sciencedirectUrl

$ws.Range("A30").Value = "IPCC JP coking coal"
$ws.Range("A31").Value = "IPCC JP PCI coal"
$ws.Range("A32").Value = "IPCC JP steam coal"
$ws.Range("A33").Value = "IPCC JP coke"

$ws.Range("Q30").Value = $ipccRietiSource
$ws.Range("Q31").Value = $ipccRietiSource
$ws.Range("Q32").Value = $ipccRietiSource
$ws.Range("Q33").Value = $ipccRietiSource

$ws.Range("Q25").Value = $ipccDensitySource

$ws.Range("Q27").Value = $ipccStationarySource
$ws.Range("Q28").Value = $ipccStationarySource
$ws.Range("Q29").Value = $ipccStationarySource

$ws.Range("A34").Value = "IPCC coking coal"

# --- Numeric values & formulas (do not consume shared-string slots) --------

# Row 20: JP coking coal
$ws.Range("B20").Value = 31.7
$ws.Range("C20").Value = 30.2
$ws.Range("D20").Formula = "=0.094*C20"

# Row 21: JP steam coal
$ws.Range("B21").Value = 27.1
$ws.Range("C21").Value = 25.8
$ws.Range("D21").Formula = "=0.096*C21"

# Row 22: JP waste plastics (no new data)

# Row 23: JP coke
$ws.Range("B23").Value = 30
$ws.Range("C23").Value = 29.8
$ws.Range("D23").Formula = "=0.109*C23"

# Row 24: RU hard coal already had data, untouched

# Row 25: RU natural gas
$ws.Range("B25").Formula = "=40.36*(1/0.554)"
$ws.Range("C25").Formula = "=36.4/0.7"
$ws.Range("D25").Formula = "=55.2/C25"

# Row 26: IPCC coke, already had data, untouched

# Row 27: US coking coal
$ws.Range("C27").Formula = "=13500/430"
$ws.Range("D27").Formula = "=94.6/C27"

# Row 28: US coke
$ws.Range("C28").Formula = "=13000/430"
$ws.Range("D28").Formula = "=C28*29.2*(44/12)/1000"

# Row 29: US natural gas
$ws.Range("C29").Value = 47.1
$ws.Range("D29").Formula = "=56.1/C29"

# Row 30: IPCC JP coking coal
$ws.Range("B30").Value = 28.94
$ws.Range("C30").Value = 26.68
$ws.Range("D30").Formula = "=(((44/12)*26.5)*C30)/1000"

# Row 31: IPCC JP PCI coal
$ws.Range("B31").Value = 28.01
$ws.Range("C31").Value = 25.74
$ws.Range("D31").Formula = "=(((44/12)*27.27)*C31)/1000"

# Row 32: IPCC JP steam coal
$ws.Range("B32").Value = 25.97
$ws.Range("C32").Value = 24.66
$ws.Range("D32").Formula = "=(((44/12)*25.68)*C32)/1000"

# Row 33: IPCC JP coke
$ws.Range("B33").Value = 29.18
$ws.Range("C33").Value = 28.81
$ws.Range("D33").Formula = "=(((44/12)*30.6)*C33)/1000"

# Row 34: IPCC coking coal
$ws.Range("C34").Value = 28.2
$ws.Range("D34").Formula = "=94.6/C34"
$ws.Range("Q34").Value = "IPCC EFDB"

# --- Source-note comments on the new CO2-factor formula cells ------------------
foreach ($cellRef in @("D30", "D31", "D32", "D33")) {
    $cmt = $ws.Range($cellRef).AddComment($noteText)
    $fnt = $cmt.Shape.TextFrame.Characters().Font
    $fnt.Name = "Tahoma"
    $fnt.Size = 9
}
